$wb = $excel.ActiveWorkbook

# Rename the "wt" and "dcin5" worksheets to reflect log2 expression data
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Make the renamed "wt_log2_expression" sheet the active sheet/selection
$wsWt.Activate()
$wsWt.Range("G22").Select()
